$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 292.21054
$ws.Range("I28").Value = 223.5
$ws.Range("J28").Value = 484.6
$ws.Range("K28").Value = 223.5
$ws.Range("L28").Value = 484.6
$ws.Range("M28").Value = 261.5
$ws.Range("N28").Value = -1454.6
# Row 32
$ws.Range("H32").Value = 882.34784
$ws.Range("I32").Value = 1138.8
$ws.Range("J32").Value = 811.1111
$ws.Range("K32").Value = 1138.8
$ws.Range("L32").Value = 811.1111
$ws.Range("M32").Value = -812.8
$ws.Range("N32").Value = -1463.1111
# Row 51
$ws.Range("H51").Value = 11537.7
$ws.Range("I51").Value = 1338.5
$ws.Range("J51").Value = 14087.5
$ws.Range("K51").Value = 1338.5
$ws.Range("L51").Value = 14087.5
$ws.Range("M51").Value = -854.5
$ws.Range("N51").Value = -15055.5
# Row 62
$ws.Range("H62").Value = 11061.091
$ws.Range("I62").Value = 12161.071
$ws.Range("J62").Value = 4901.2
$ws.Range("K62").Value = 12161.071
$ws.Range("L62").Value = 4901.2
$ws.Range("M62").Value = -11537.071
$ws.Range("N62").Value = -6149.2
# Row 65
$ws.Range("H65").Value = 11061.091
$ws.Range("I65").Value = 12161.071
$ws.Range("J65").Value = 4901.2
$ws.Range("K65").Value = 60805.355
$ws.Range("L65").Value = 24506
$ws.Range("M65").Value = -57685.355
$ws.Range("N65").Value = -30746
# Row 98
$ws.Range("H98").Value = 915.11536
$ws.Range("I98").Value = 915.11536
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 915.11536
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 582.88464
# Row 122
$ws.Range("H122").Value = 915.11536
$ws.Range("I122").Value = 915.11536
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2745.34608
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -295.3460800000003
# Row 129
$ws.Range("H129").Value = 635.78125
$ws.Range("I129").Value = 493.26086
$ws.Range("K129").Value = 1479.78258
$ws.Range("M129").Value = 3520.21742
# Row 135
$ws.Range("H135").Value = 16130556
$ws.Range("I135").Value = 494.23077
$ws.Range("J135").Value = 27780044
$ws.Range("K135").Value = 4448.07693
$ws.Range("L135").Value = 250020396
$ws.Range("M135").Value = -1913.07693
$ws.Range("N135").Value = -250025466

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6106.421
$ws.Range("I2").Value = 1000.2308
$ws.Range("J2").Value = 17169.834
$ws.Range("K2").Value = 1000.2308
$ws.Range("L2").Value = 17169.834
$ws.Range("M2").Value = -887.2308
$ws.Range("N2").Value = -17395.834
# Row 45
$ws.Range("H45").Value = 1032.1
$ws.Range("I45").Value = 1219.3334
$ws.Range("J45").Value = 878.9091
$ws.Range("K45").Value = 1219.3334
$ws.Range("L45").Value = 878.9091
$ws.Range("M45").Value = -842.3334
$ws.Range("N45").Value = -1632.9091
# Row 52
$ws.Range("H52").Value = 45780
$ws.Range("J52").Value = 45780
$ws.Range("L52").Value = 45780
$ws.Range("N52").Value = -46416
# Row 76
$ws.Range("H76").Value = 40079.332
$ws.Range("J76").Value = 40079.332
$ws.Range("L76").Value = 40079.332
$ws.Range("N76").Value = -40755.332
# Row 79
$ws.Range("H79").Value = 40079.332
$ws.Range("J79").Value = 40079.332
$ws.Range("L79").Value = 40079.332
$ws.Range("N79").Value = -42419.332
# Row 102
$ws.Range("H102").Value = 1640
$ws.Range("I102").Value = 1640
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1640
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -18
# Row 116
$ws.Range("H116").Value = 6106.421
$ws.Range("I116").Value = 1000.2308
$ws.Range("J116").Value = 17169.834
$ws.Range("K116").Value = 1000.2308
$ws.Range("L116").Value = 17169.834
$ws.Range("M116").Value = 1293.7692
$ws.Range("N116").Value = -21757.834
# Row 122
$ws.Range("H122").Value = 1419.9
$ws.Range("I122").Value = 1374.875
$ws.Range("K122").Value = 4124.625
$ws.Range("M122").Value = -1674.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6106.421
$ws.Range("I3").Value = 1000.2308
$ws.Range("J3").Value = 17169.834
$ws.Range("K3").Value = 1000.2308
$ws.Range("L3").Value = 17169.834
$ws.Range("M3").Value = -886.2308
$ws.Range("N3").Value = -17397.834
# Row 134
$ws.Range("H134").Value = 21233.418
$ws.Range("J134").Value = 93636.836
$ws.Range("L134").Value = 280910.508
$ws.Range("N134").Value = -285980.508

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 30202.2
$ws.Range("I3").Value = 3001
$ws.Range("J3").Value = 48336.332
$ws.Range("K3").Value = 3001
$ws.Range("L3").Value = 48336.332
$ws.Range("M3").Value = -2888
$ws.Range("N3").Value = -48562.332
# Row 16
$ws.Range("H16").Value = 1048.3636
$ws.Range("I16").Value = 1003.1667
$ws.Range("K16").Value = 1003.1667
$ws.Range("M16").Value = -716.1667
# Row 31
$ws.Range("H31").Value = 16024.893
$ws.Range("I31").Value = 18340.928
$ws.Range("J31").Value = 13708.857
$ws.Range("K31").Value = 18340.928
$ws.Range("L31").Value = 13708.857
$ws.Range("M31").Value = -18045.928
$ws.Range("N31").Value = -14298.857
# Row 34
$ws.Range("H34").Value = 16024.893
$ws.Range("I34").Value = 18340.928
$ws.Range("J34").Value = 13708.857
$ws.Range("K34").Value = 18340.928
$ws.Range("L34").Value = 13708.857
$ws.Range("M34").Value = -18138.928
$ws.Range("N34").Value = -14112.857
# Row 105
$ws.Range("H105").Value = 995
$ws.Range("I105").Value = 995
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 995
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = 752
# Row 107
$ws.Range("H107").Value = 393.82608
$ws.Range("I107").Value = 252.9375
$ws.Range("J107").Value = 715.8570999999999
$ws.Range("K107").Value = 252.9375
$ws.Range("L107").Value = 715.8570999999999
$ws.Range("M107").Value = 1667.0625
$ws.Range("N107").Value = -4555.8571
# Row 113
$ws.Range("H113").Value = 1048.3636
$ws.Range("I113").Value = 1003.1667
$ws.Range("K113").Value = 1003.1667
$ws.Range("M113").Value = 1166.8333
# Row 122
$ws.Range("H122").Value = 1640
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 1625
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 4875
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -9775
# Row 134
$ws.Range("H134").Value = 1299.3334
$ws.Range("I134").Value = 1236.1111
$ws.Range("J134").Value = 1489
$ws.Range("K134").Value = 3708.3333
$ws.Range("L134").Value = 4467
$ws.Range("M134").Value = -1173.3333
$ws.Range("N134").Value = -9537

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 908.0526
$ws.Range("I97").Value = 880.2727
$ws.Range("J97").Value = 946.25
$ws.Range("K97").Value = 880.2727
$ws.Range("L97").Value = 946.25
$ws.Range("M97").Value = -384.2727
$ws.Range("N97").Value = -1938.25
# Row 113
$ws.Range("H113").Value = 1119.6666
$ws.Range("I113").Value = 1213.6
$ws.Range("J113").Value = 650
$ws.Range("K113").Value = 1213.6
$ws.Range("L113").Value = 650
$ws.Range("M113").Value = 956.4000000000001
$ws.Range("N113").Value = -4990
# Row 122
$ws.Range("H122").Value = 1452
$ws.Range("J122").Value = 1602.6666
$ws.Range("L122").Value = 4807.9998
$ws.Range("N122").Value = -9707.9998
# Row 126
$ws.Range("H126").Value = 1427.0889
$ws.Range("I126").Value = 1303.0938
$ws.Range("J126").Value = 1732.3077
$ws.Range("K126").Value = 3909.2814
$ws.Range("L126").Value = 5196.9231
$ws.Range("M126").Value = -1439.2814
$ws.Range("N126").Value = -10136.9231
# Row 132
$ws.Range("H132").Value = 20945.941
$ws.Range("I132").Value = 1188.6666
$ws.Range("J132").Value = 68363.39999999999
$ws.Range("K132").Value = 3565.9998
$ws.Range("L132").Value = 205090.2
$ws.Range("M132").Value = -1035.9998
$ws.Range("N132").Value = -210150.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 36502.75
$ws.Range("I10").Value = 2003
$ws.Range("J10").Value = 48002.668
$ws.Range("K10").Value = 2003
$ws.Range("L10").Value = 48002.668
$ws.Range("M10").Value = -1863
$ws.Range("N10").Value = -48282.668
# Row 93
$ws.Range("H93").Value = 1692.9286
$ws.Range("I93").Value = 1618.6316
$ws.Range("J93").Value = 1754.3043
$ws.Range("K93").Value = 1618.6316
$ws.Range("L93").Value = 1754.3043
$ws.Range("M93").Value = -370.6315999999999
$ws.Range("N93").Value = -4250.3043
# Row 122
$ws.Range("H122").Value = 2771.9429
$ws.Range("I122").Value = 2809.6553
$ws.Range("J122").Value = 2589.6667
$ws.Range("K122").Value = 8428.965899999999
$ws.Range("L122").Value = 7769.000100000001
$ws.Range("M122").Value = -5978.965899999999
$ws.Range("N122").Value = -12669.0001
# Row 132
$ws.Range("H132").Value = 347170.12
$ws.Range("I132").Value = 108385.79
$ws.Range("J132").Value = 630726.5
$ws.Range("K132").Value = 325157.37
$ws.Range("L132").Value = 1892179.5
$ws.Range("M132").Value = -322627.37
$ws.Range("N132").Value = -1897239.5
# Row 136
$ws.Range("H136").Value = 197222.92
$ws.Range("I136").Value = 286453.84
$ws.Range("J136").Value = 2030.3125
$ws.Range("K136").Value = 859361.52
$ws.Range("L136").Value = 6090.9375
$ws.Range("M136").Value = -856811.52
$ws.Range("N136").Value = -11190.9375

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 31798
$ws.Range("J82").Value = 31798
$ws.Range("L82").Value = 31798
$ws.Range("N82").Value = -32564
# Row 85
$ws.Range("H85").Value = 31798
$ws.Range("J85").Value = 31798
$ws.Range("L85").Value = 31798
$ws.Range("N85").Value = -34450
# Row 113
$ws.Range("H113").Value = 331.5
$ws.Range("I113").Value = 360.4
$ws.Range("J113").Value = 283.33334
$ws.Range("K113").Value = 1081.2
$ws.Range("L113").Value = 850.0000200000001
$ws.Range("M113").Value = 1088.8
$ws.Range("N113").Value = -5190.00002
# Row 122
$ws.Range("H122").Value = 4231.8965
$ws.Range("I122").Value = 2151.7856
$ws.Range("J122").Value = 6173.3335
$ws.Range("K122").Value = 6455.3568
$ws.Range("L122").Value = 18520.0005
$ws.Range("M122").Value = -4005.3568
$ws.Range("N122").Value = -23420.0005
# Row 126
$ws.Range("H126").Value = 798.3043
$ws.Range("I126").Value = 667.85
$ws.Range("J126").Value = 1668
$ws.Range("K126").Value = 2003.55
$ws.Range("L126").Value = 5004
$ws.Range("M126").Value = 466.4499999999998
$ws.Range("N126").Value = -9944
# Row 132
$ws.Range("H132").Value = 5196.577
$ws.Range("I132").Value = 1516.2307
$ws.Range("J132").Value = 8876.923000000001
$ws.Range("K132").Value = 4548.6921
$ws.Range("L132").Value = 26630.769
$ws.Range("M132").Value = -2018.6921
$ws.Range("N132").Value = -31690.769
